# Generate Report for Handoff
# Adds two new localized-file rows (3406dc6f-... and f2071dd3-...) to every
# sheet of the localization-status workbook: the "Overview" summary sheet
# plus the "zh-cn" and "de-de" per-locale sheets. Each sheet's backing
# Excel Table is resized to include the two new rows, and the new
# hyperlinked filename cells get a Hyperlinks.Add() entry just like the
# existing rows. Date/time columns keep the same custom
# "yyyy-mm-dd HH:mm:ss" number format used by the existing rows.
#
# A leading "'" forces Excel to store the literal text "True"/"False"/""
# instead of auto-coercing it to a Boolean (or dropping an empty value
# altogether) - it keeps the same text content, matching the existing
# rows' shared-string cells.

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (A1:G3 -> A1:G5)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G5"))

$ws.Cells.Item(4,1).Value = "3406dc6f-a315-40ae-87ad-281da0299a7f.md"
$ws.Hyperlinks.Add($ws.Cells.Item(4,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/3406dc6f-a315-40ae-87ad-281da0299a7f.md", "", "", "e2e\3406dc6f-a315-40ae-87ad-281da0299a7f.md")
$ws.Cells.Item(4,3).Value = ".md"
$ws.Cells.Item(4,4).Value = "'"
$ws.Cells.Item(4,5).Value = "Ready for handoff"
$ws.Cells.Item(4,6).Value = "Ready for handoff"
$ws.Cells.Item(4,7).Value = "2016-08-17 06:36:12"
$ws.Cells.Item(4,7).NumberFormat = $dateFmt

$ws.Cells.Item(5,1).Value = "f2071dd3-b806-42c1-be05-f7a26108fdf5.md"
$ws.Hyperlinks.Add($ws.Cells.Item(5,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/f2071dd3-b806-42c1-be05-f7a26108fdf5.md", "", "", "e2e\f2071dd3-b806-42c1-be05-f7a26108fdf5.md")
$ws.Cells.Item(5,3).Value = ".md"
$ws.Cells.Item(5,4).Value = "'"
$ws.Cells.Item(5,5).Value = "Ready for handoff"
$ws.Cells.Item(5,6).Value = "Ready for handoff"
$ws.Cells.Item(5,7).Value = "2016-08-17 06:36:12"
$ws.Cells.Item(5,7).NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Sheet "zh-cn" (A1:P3 -> A1:P5)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P5"))

$ws.Hyperlinks.Add($ws.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/3406dc6f-a315-40ae-87ad-281da0299a7f.md", "", "", "3406dc6f-a315-40ae-87ad-281da0299a7f.md")
$ws.Cells.Item(4,2).Value = ".md"
$ws.Cells.Item(4,3).Value = "Ready for handoff"
$ws.Cells.Item(4,4).Value = "e2e"
$ws.Cells.Item(4,5).Value = "ht"
$ws.Cells.Item(4,6).Value = "'False"
$ws.Cells.Item(4,7).Value = "3406dc6f-a315-40ae-87ad-281da0299a7f.6380b91eacdfa381f25b0779db83fed5ace595ba.zh-cn.xlf"
$ws.Cells.Item(4,8).Value = "2016-08-17 06:36:04"
$ws.Cells.Item(4,8).NumberFormat = $dateFmt
$ws.Cells.Item(4,9).Value = "'"
$ws.Cells.Item(4,10).Value = "'"
$ws.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4,11).NumberFormat = $dateFmt
$ws.Cells.Item(4,12).Value = "'"
$ws.Cells.Item(4,13).Value = "'True"
$ws.Cells.Item(4,14).Value = "'"
$ws.Cells.Item(4,15).Value = "'False"
$ws.Cells.Item(4,16).Value = "'"

$ws.Hyperlinks.Add($ws.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/f2071dd3-b806-42c1-be05-f7a26108fdf5.md", "", "", "f2071dd3-b806-42c1-be05-f7a26108fdf5.md")
$ws.Cells.Item(5,2).Value = ".md"
$ws.Cells.Item(5,3).Value = "Ready for handoff"
$ws.Cells.Item(5,4).Value = "e2e"
$ws.Cells.Item(5,5).Value = "ht"
$ws.Cells.Item(5,6).Value = "'False"
$ws.Cells.Item(5,7).Value = "f2071dd3-b806-42c1-be05-f7a26108fdf5.270134bc23556e116fb1583be525b79973b81444.zh-cn.xlf"
$ws.Cells.Item(5,8).Value = "2016-08-17 06:36:04"
$ws.Cells.Item(5,8).NumberFormat = $dateFmt
$ws.Cells.Item(5,9).Value = "'"
$ws.Cells.Item(5,10).Value = "'"
$ws.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(5,11).NumberFormat = $dateFmt
$ws.Cells.Item(5,12).Value = "'"
$ws.Cells.Item(5,13).Value = "'True"
$ws.Cells.Item(5,14).Value = "'"
$ws.Cells.Item(5,15).Value = "'False"
$ws.Cells.Item(5,16).Value = "'"

# ---------------------------------------------------------------------
# Sheet "de-de" (A1:P3 -> A1:P5)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P5"))

$ws.Hyperlinks.Add($ws.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/3406dc6f-a315-40ae-87ad-281da0299a7f.md", "", "", "3406dc6f-a315-40ae-87ad-281da0299a7f.md")
$ws.Cells.Item(4,2).Value = ".md"
$ws.Cells.Item(4,3).Value = "Ready for handoff"
$ws.Cells.Item(4,4).Value = "e2e"
$ws.Cells.Item(4,5).Value = "ht"
$ws.Cells.Item(4,6).Value = "'False"
$ws.Cells.Item(4,7).Value = "3406dc6f-a315-40ae-87ad-281da0299a7f.6380b91eacdfa381f25b0779db83fed5ace595ba.de-de.xlf"
$ws.Cells.Item(4,8).Value = "2016-08-17 06:36:12"
$ws.Cells.Item(4,8).NumberFormat = $dateFmt
$ws.Cells.Item(4,9).Value = "'"
$ws.Cells.Item(4,10).Value = "'"
$ws.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4,11).NumberFormat = $dateFmt
$ws.Cells.Item(4,12).Value = "'"
$ws.Cells.Item(4,13).Value = "'True"
$ws.Cells.Item(4,14).Value = "'"
$ws.Cells.Item(4,15).Value = "'False"
$ws.Cells.Item(4,16).Value = "'"

$ws.Hyperlinks.Add($ws.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/f2071dd3-b806-42c1-be05-f7a26108fdf5.md", "", "", "f2071dd3-b806-42c1-be05-f7a26108fdf5.md")
$ws.Cells.Item(5,2).Value = ".md"
$ws.Cells.Item(5,3).Value = "Ready for handoff"
$ws.Cells.Item(5,4).Value = "e2e"
$ws.Cells.Item(5,5).Value = "ht"
$ws.Cells.Item(5,6).Value = "'False"
$ws.Cells.Item(5,7).Value = "f2071dd3-b806-42c1-be05-f7a26108fdf5.270134bc23556e116fb1583be525b79973b81444.de-de.xlf"
$ws.Cells.Item(5,8).Value = "2016-08-17 06:36:12"
$ws.Cells.Item(5,8).NumberFormat = $dateFmt
$ws.Cells.Item(5,9).Value = "'"
$ws.Cells.Item(5,10).Value = "'"
$ws.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(5,11).NumberFormat = $dateFmt
$ws.Cells.Item(5,12).Value = "'"
$ws.Cells.Item(5,13).Value = "'True"
$ws.Cells.Item(5,14).Value = "'"
$ws.Cells.Item(5,15).Value = "'False"
$ws.Cells.Item(5,16).Value = "'"
